$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new "To Do" rows for CRUD work items:
#   - row 3 : Branch CRUD (update, delete)
#   - row 5 : Commit CRUD (update message, delete)   (after the insert above, old row 3 becomes row 4)
#   - row 9 : Tag CRUD (update name, delete)
# Each insert pushes everything below it down by one row, so doing them
# top-to-bottom at these row numbers lands the new rows exactly where the
# diff expects them.
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(9).Insert()

# Populate the newly inserted rows.
$ws.Range("A3").Value = "Branch CRUD (update, delete)"
$ws.Range("B3").Value = "To Do"
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "TBD"
$ws.Range("E3").Value = "Edit branch name; delete branch"
$ws.Range("A3").Font.Bold = $true

$ws.Range("A5").Value = "Commit CRUD (update message, delete)"
$ws.Range("B5").Value = "To Do"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "TBD"
$ws.Range("E5").Value = "Edit commit message; delete commit"
$ws.Range("A5").Font.Bold = $true

$ws.Range("A9").Value = "Tag CRUD (update name, delete)"
$ws.Range("B9").Value = "To Do"
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "TBD"
$ws.Range("E9").Value = "Rename or remove tags"
$ws.Range("A9").Font.Bold = $true

# Remove the trailing "Docs / Deployment / Promotion" rows, which (after the
# three inserts above) now live at rows 21-23.
$ws.Range("A21:E23").EntireRow.Delete()

# Update the saved selection to match the edited workbook's view.
$ws.Range("D27").Select()
